$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write literal text into a cell without Excel's "smart" number/date/
# percent auto-conversion re-typing the text would trigger. We build the text
# via a formula (so it is never parsed as a literal token), then immediately
# convert that formula to a plain value in place (Copy + PasteSpecial values),
# which preserves the cell's existing style/number-format.
function Set-TextValue($addr, $text) {
    $escaped = $text -replace '"', '""'
    $ws.Range($addr).Formula = '="' + $escaped + '"'
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

# --- Row 2: search input + match score header cell ---
Set-TextValue 'A2' 'Vladimir Poutine'
Set-TextValue 'F2' '100 %'

# --- Existing rows 3-5: refreshed match results ---
# Row 3
Set-TextValue 'B3' '0. (100%) -  Vladimir Vladimirovich Poutine'
Set-TextValue 'C3' 'FR - Liste de sanctions de la Direction Générale du Trésor (DGT)'
Set-TextValue 'D3' '07/10/1952'
$ws.Range('E3').Value = ""
Set-TextValue 'G3' 'todoByFrontDev/713310123456789012345678'

# Row 4
Set-TextValue 'B4' '1. (100%) -  Vladimir Vladimirovich Putin'
Set-TextValue 'C4' 'UE - Liste consolidée des personnes, groupes et entités faisant l''objet de sanctions financières de l''UE'
Set-TextValue 'D4' '07/10/1952'
$ws.Range('E4').Value = ""
Set-TextValue 'G4' 'todoByFrontDev/171350123456789012345678'

# Row 5
Set-TextValue 'B5' '2. (91.36%) -  Vladimir Vladimirovich Putin'
Set-TextValue 'C5' 'US - Liste OFAC des ressortissants spécialement désignés et des personnes bloquées (SDN)'
Set-TextValue 'D5' '07/10/1952'
Set-TextValue 'E5' 'Russia'
Set-TextValue 'G5' 'todoByFrontDev/975901234567890123456789'

# --- New rows 6-17: additional match results (copy formatting down first) ---
# Row 6
$ws.Range('A5:H5').Copy($ws.Range('A6:H6'))
Set-TextValue 'B6' '3. (76.46%) -  Vladimir Putin'
Set-TextValue 'C6' 'GB - Liste consolidée des sanctions financières du Royaume-Uni (HMT)'
Set-TextValue 'D6' '07/10/1952'
Set-TextValue 'E6' 'Russia'
Set-TextValue 'G6' 'todoByFrontDev/655180123456789012345678'
$ws.Range('H6').Value = 3

# Row 7
$ws.Range('A6:H6').Copy($ws.Range('A7:H7'))
Set-TextValue 'B7' '4. (76.46%) -  Vladimir Putin'
Set-TextValue 'C7' 'GB - Liste consolidée des sanctions financières du Royaume-Uni (HMT)'
Set-TextValue 'D7' '07/10/1952'
Set-TextValue 'E7' 'Russia'
Set-TextValue 'G7' 'todoByFrontDev/655170123456789012345678'
$ws.Range('H7').Value = 3

# Row 8
$ws.Range('A7:H7').Copy($ws.Range('A8:H8'))
Set-TextValue 'B8' '5. (70.04%) -  Vladimir Vladimirovich Putin'
Set-TextValue 'C8' 'CA - Liste consolidée des sanctions autonomes canadiennes'
Set-TextValue 'D8' '1952'
Set-TextValue 'E8' 'Russia'
Set-TextValue 'G8' 'todoByFrontDev/742020123456789012345678'
$ws.Range('H8').Value = 3

# Row 9
$ws.Range('A8:H8').Copy($ws.Range('A9:H9'))
Set-TextValue 'B9' '6. (67.48%) -  Vladimir Vladimirovich  Vladimirov   '
Set-TextValue 'C9' 'CA - Liste consolidée des sanctions autonomes canadiennes'
Set-TextValue 'D9' '14/10/1975'
Set-TextValue 'E9' 'Russia'
Set-TextValue 'G9' 'todoByFrontDev/746210123456789012345678'
$ws.Range('H9').Value = 3

# Row 10
$ws.Range('A9:H9').Copy($ws.Range('A10:H10'))
Set-TextValue 'B10' '7. (66.33%) -  Vladimir Vladimirovich Putin'
Set-TextValue 'C10' 'GB - Liste consolidée des sanctions financières du Royaume-Uni (HMT)'
Set-TextValue 'D10' '07/10/1952'
Set-TextValue 'E10' 'Russia'
Set-TextValue 'G10' 'todoByFrontDev/655190123456789012345678'
$ws.Range('H10').Value = 3

# Row 11
$ws.Range('A10:H10').Copy($ws.Range('A11:H11'))
Set-TextValue 'B11' '8. (60.26%) -  Vladimir Vladimirovich Vladimirov'
Set-TextValue 'C11' 'GB - Liste consolidée des sanctions financières du Royaume-Uni (HMT)'
Set-TextValue 'D11' '14/10/1975'
Set-TextValue 'E11' 'Russia'
Set-TextValue 'G11' 'todoByFrontDev/685770123456789012345678'
$ws.Range('H11').Value = 3

# Row 12
$ws.Range('A11:H11').Copy($ws.Range('A12:H12'))
Set-TextValue 'B12' '9. (54.38%) -  Екатерина Владимировна Тихонова'
Set-TextValue 'C12' 'UE - Liste consolidée des personnes, groupes et entités faisant l''objet de sanctions financières de l''UE'
Set-TextValue 'D12' '31/08/1986'
Set-TextValue 'E12' 'Russia'
Set-TextValue 'G12' 'todoByFrontDev/173910123456789012345678'
$ws.Range('H12').Value = 3

# Row 13
$ws.Range('A12:H12').Copy($ws.Range('A13:H13'))
Set-TextValue 'B13' '10. (54.38%) -  Mariya Vorontsova'
Set-TextValue 'C13' 'UE - Liste consolidée des personnes, groupes et entités faisant l''objet de sanctions financières de l''UE'
Set-TextValue 'D13' '28/04/1985'
Set-TextValue 'E13' 'Russia'
Set-TextValue 'G13' 'todoByFrontDev/173830123456789012345678'
$ws.Range('H13').Value = 3

# Row 14
$ws.Range('A13:H13').Copy($ws.Range('A14:H14'))
Set-TextValue 'B14' '11. (53.79%) -  Lazher Ben Khalifa Ben Ahmed Rouine'
Set-TextValue 'C14' 'US - Liste OFAC des ressortissants spécialement désignés et des personnes bloquées (SDN)'
Set-TextValue 'D14' '20/11/1975'
Set-TextValue 'E14' 'Tunisia'
Set-TextValue 'G14' 'todoByFrontDev/117101234567890123456789'
$ws.Range('H14').Value = 3

# Row 15
$ws.Range('A14:H14').Copy($ws.Range('A15:H15'))
Set-TextValue 'B15' '12. (53.79%) -  Maria Vladimirovna Vorontsova'
Set-TextValue 'C15' 'US - Liste OFAC des ressortissants spécialement désignés et des personnes bloquées (SDN)'
Set-TextValue 'D15' '28/04/1985'
Set-TextValue 'E15' 'Russia'
Set-TextValue 'G15' 'todoByFrontDev/104600123456789012345678'
$ws.Range('H15').Value = 3

# Row 16
$ws.Range('A15:H15').Copy($ws.Range('A16:H16'))
Set-TextValue 'B16' '13. (53.79%) -  Katerina Vladimirovna Tikhonova'
Set-TextValue 'C16' 'US - Liste OFAC des ressortissants spécialement désignés et des personnes bloquées (SDN)'
Set-TextValue 'D16' '31/08/1986'
Set-TextValue 'E16' 'Russia'
Set-TextValue 'G16' 'todoByFrontDev/103840123456789012345678'
$ws.Range('H16').Value = 3

# Row 17
$ws.Range('A16:H16').Copy($ws.Range('A17:H17'))
Set-TextValue 'B17' '14. (53.79%) -  Al-azhar Ben Khalifa Ben Ahmed Rouine'
Set-TextValue 'C17' 'UE - Liste consolidée des personnes, groupes et entités faisant l''objet de sanctions financières de l''UE'
Set-TextValue 'D17' '20/11/1975'
Set-TextValue 'E17' 'Russia'
Set-TextValue 'G17' 'todoByFrontDev/145210123456789012345678'
$ws.Range('H17').Value = 3

